$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(101, 9, 30, 15, 60, 15),
    @(902, 1, 0, 0, 0, 0),
    @(1201, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(401, 9, 48, 67, 75, 45),
    @(701, 3, 90, 45, 97, 15),
    @(1202, 2, 10, 10, 10, 10),
    @(901, 16, 15, 45, 60, 60),
    @(301, 6, 45, 30, 60, 45),
    @(801, 3, 67, 65, 52, 45),
    @(501, 9, 52, 30, 75, 45),
    @(601, 9, 60, 67, 60, 42),
    @(201, 9, 30, 15, 45, 30),
    @(502, 0, 4, 0, 0, 0),
    @(1101, 0, 15, 30, 30, 0),
    @(802, 0, 4, 5, 4, 0),
    @(3, 0, 3, 3, 3, 3),
    @(1, 0, 2, 2, 2, 2),
    @(2, 0, 2, 2, 2, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
